# Update countries & provincias Spain
#
# Refreshes the "Pais" COVID-19 dashboard with a newer data pull (the report
# timestamp moves from 14:35 to 15:05) and re-sorts a handful of rows whose
# "Casos totales" (column B) changed enough to swap places with their
# neighbour in the descending sort: Arabia Saudita/Chile, Serbia/Corea del
# Sur, Tayikistan/Uzbekistan/Guinea, Sierra Leona/Principado de
# Andorra/Nicaragua, Belice/Nueva Caledonia/Santa Lucia, Islas Turcas y
# Caicos/Groenlandia and Papua Nueva Guinea/Islas Virgenes Britanicas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report timestamp
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 27 de Mayo de 2020 a las 15:05"

# Estados Unidos (row 4) - updated totals
$ws.Cells.Item(4, 2).Value = 1727700
$ws.Cells.Item(4, 3).Value = 2425
$ws.Cells.Item(4, 5).Value = 1147102

# Alemania (row 11) - updated totals
$ws.Cells.Item(11, 2).Value = 181333
$ws.Cells.Item(11, 3).Value = 45
$ws.Cells.Item(11, 5).Value = 10035

# India (row 13) - updated totals
$ws.Cells.Item(13, 2).Value = 153304
$ws.Cells.Item(13, 3).Value = 2511
$ws.Cells.Item(13, 4).Value = 64827
$ws.Cells.Item(13, 5).Value = 84110
$ws.Cells.Item(13, 7).Value = 23
$ws.Cells.Item(13, 8).Value = 4367

# Rows 18-19: Arabia Saudita overtakes Chile
$ws.Cells.Item(18, 1).Value = "Arabia Saudita"
$ws.Cells.Item(18, 2).Value = 78541
$ws.Cells.Item(18, 3).Value = 1815
$ws.Cells.Item(18, 4).Value = 51022
$ws.Cells.Item(18, 5).Value = 27094
$ws.Cells.Item(18, 7).Value = 14
$ws.Cells.Item(18, 8).Value = 425

$ws.Cells.Item(19, 1).Value = "Chile"
$ws.Cells.Item(19, 2).Value = 77961
$ws.Cells.Item(19, 4).Value = 30915
$ws.Cells.Item(19, 5).Value = 46240
$ws.Cells.Item(19, 8).Value = 806

# Rows 51-52: Serbia overtakes Corea del Sur
$ws.Cells.Item(51, 1).Value = "Serbia"
$ws.Cells.Item(51, 2).Value = 11275
$ws.Cells.Item(51, 3).Value = 48
$ws.Cells.Item(51, 4).Value = 6277
$ws.Cells.Item(51, 5).Value = 4758
$ws.Cells.Item(51, 7).Value = 1
$ws.Cells.Item(51, 8).Value = 240

$ws.Cells.Item(52, 1).Value = "Corea del Sur"
$ws.Cells.Item(52, 2).Value = 11265
$ws.Cells.Item(52, 3).Value = 40
$ws.Cells.Item(52, 4).Value = 10295
$ws.Cells.Item(52, 5).Value = 701
$ws.Cells.Item(52, 8).Value = 269

# Rows 76-78: Tayikistan moves ahead of Uzbekistan and Guinea
$ws.Cells.Item(76, 1).Value = "Tayikistan"
$ws.Cells.Item(76, 2).Value = 3424
$ws.Cells.Item(76, 3).Value = 158
$ws.Cells.Item(76, 4).Value = 1575
$ws.Cells.Item(76, 5).Value = 1802
$ws.Cells.Item(76, 8).Value = 47

$ws.Cells.Item(77, 1).Value = "Uzbekistan"
$ws.Cells.Item(77, 2).Value = 3355
$ws.Cells.Item(77, 3).Value = 65
$ws.Cells.Item(77, 4).Value = 2659
$ws.Cells.Item(77, 5).Value = 682
$ws.Cells.Item(77, 8).Value = 14

$ws.Cells.Item(78, 1).Value = "Guinea"
$ws.Cells.Item(78, 2).Value = 3275
$ws.Cells.Item(78, 4).Value = 1673
$ws.Cells.Item(78, 5).Value = 1582
$ws.Cells.Item(78, 8).Value = 20

# Row 93: small correction, no reordering
$ws.Cells.Item(93, 2).Value = 1805
$ws.Cells.Item(93, 3).Value = 1
$ws.Cells.Item(93, 5).Value = 3

# Rows 123-125: Sierra Leona moves ahead of Principado de Andorra and Nicaragua
$ws.Cells.Item(123, 1).Value = "Sierra Leona"
$ws.Cells.Item(123, 2).Value = 782
$ws.Cells.Item(123, 3).Value = 28
$ws.Cells.Item(123, 4).Value = 297
$ws.Cells.Item(123, 5).Value = 440
$ws.Cells.Item(123, 7).Value = 1
$ws.Cells.Item(123, 8).Value = 45

$ws.Cells.Item(124, 1).Value = "Principado de Andorra"
$ws.Cells.Item(124, 2).Value = 763
$ws.Cells.Item(124, 4).Value = 676
$ws.Cells.Item(124, 5).Value = 36
$ws.Cells.Item(124, 8).Value = 51

$ws.Cells.Item(125, 1).Value = "Nicaragua"
$ws.Cells.Item(125, 2).Value = 759
$ws.Cells.Item(125, 4).Value = 370
$ws.Cells.Item(125, 5).Value = 354
$ws.Cells.Item(125, 8).Value = 35

# Rows 199-201: Belice moves ahead of Nueva Caledonia and Santa Lucia
$ws.Cells.Item(199, 1).Value = "Belice"
$ws.Cells.Item(199, 4).Value = 16
$ws.Cells.Item(199, 8).Value = 2

$ws.Cells.Item(200, 1).Value = "Nueva Caledonia"

$ws.Cells.Item(201, 1).Value = "Santa Lucia"
$ws.Cells.Item(201, 4).Value = 18
$ws.Cells.Item(201, 8).Value = 0

# Rows 207-208: Islas Turcas y Caicos overtakes Groenlandia
$ws.Cells.Item(207, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(207, 4).Value = 10
$ws.Cells.Item(207, 8).Value = 1

$ws.Cells.Item(208, 1).Value = "Groenlandia"
$ws.Cells.Item(208, 4).Value = 11
$ws.Cells.Item(208, 8).Value = 0

# Rows 213-214: Papua Nueva Guinea overtakes Islas Virgenes Britanicas
$ws.Cells.Item(213, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213, 4).Value = 8
$ws.Cells.Item(213, 8).Value = 0

$ws.Cells.Item(214, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 8).Value = 1
